# Apply the changes described by the commit "Complet Arm State and Adding Item State"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item_Initialize_Data")

# --- Adding Item State: highlight a handful of item rows in yellow ---
# These B-column cells (VLOOKUP name lookups) move from the existing
# green-ish "style 6" fill to a new yellow fill (same font/border, new fill).
$itemRows = @(9, 10, 12, 14, 34, 48)
foreach ($r in $itemRows) {
    $ws.Cells.Item($r, 2).Interior.Color = 65535   # RGB(255,255,0) => yellow FFFFFF00
}

# --- Complete Arm State: fix up the C/D (type/value) columns for rows 59-67 ---
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 1

$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(60, 4).Value = 2

$ws.Cells.Item(61, 3).Value = 2
$ws.Cells.Item(61, 4).Value = 0

$ws.Cells.Item(62, 4).Value = 3

$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 3

$ws.Cells.Item(64, 4).Value = 1

$ws.Cells.Item(65, 4).Value = 2

$ws.Cells.Item(66, 4).Value = 3

$ws.Cells.Item(67, 4).Value = 4

# --- Update the saved view/selection state ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G55").Select() | Out-Null
